# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The "Periodo Mora" rows (B16:J27) get re-sorted by period (column E) in
# ascending order (2107 .. 2207), which is why the "Valor Mora" amounts in
# column F move: the 40000 value that belonged to period 2207 now shows up
# on the first row (period 2107) and the 48000 value that belonged to
# period 2107 now shows up on the last row (period 2207) - all the periods
# in between already carried 48000 so they look unchanged.
#
# Row styling (borders/shading) must stay anchored to the row position
# (row 27 keeps the special "last row" bottom border), so we update the
# cell values directly instead of using Range.Sort (which would drag the
# per-row formatting along with the data).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# New ascending order of "Periodo Mora" values for rows 16-27.
$periodos = @("2107","2109","2110","2111","2112","2201","2202","2203","2204","2205","2206","2207")

for ($i = 0; $i -lt $periodos.Length; $i++) {
    $row = 16 + $i
    $ws.Cells.Item($row, 5).Value = $periodos[$i]
}

# "Valor Mora" column: only the first and last rows actually change value
# once the data is resorted (the rest already were 48000).
$ws.Cells.Item(16, 6).Value = 48000
$ws.Cells.Item(27, 6).Value = 40000
